$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.713.93"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.004.54"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "513.10"
$ws.Range("E5").Value = "  +5.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.94"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +4.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.51"
$ws.Range("E9").Value = "  +7.19%  "
$ws.Range("E10").Value = "  +9.48%  "
$ws.Range("E11").Value = "  +3.19%  "
$ws.Range("E12").Value = "  +2.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.522.73"
$ws.Range("E13").Value = "  +2.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.83"
$ws.Range("E14").Value = "  +5.95%  "
$ws.Range("E15").Value = "  +14.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.749.22"
$ws.Range("E16").Value = "  +2.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.005.32"
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.00"
$ws.Range("E18").Value = "  +7.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.54"
$ws.Range("E19").Value = "  +5.26%  "
$ws.Range("E20").Value = "  +6.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.48"
$ws.Range("E21").Value = "  +6.63%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  +5.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.18"
$ws.Range("E24").Value = "  +5.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0909"
$ws.Range("E27").Value = "  +8.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.71"
$ws.Range("E28").Value = "  +3.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.98"
$ws.Range("E29").Value = "  +7.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.24"
$ws.Range("E30").Value = "  +7.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.81"
$ws.Range("E31").Value = "  +8.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.74"
$ws.Range("E32").Value = "  +8.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.14"
$ws.Range("E33").Value = "  +4.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.56"
$ws.Range("E34").Value = "  +5.35%  "
$ws.Range("E35").Value = "  +1.88%  "
$ws.Range("E36").Value = "  +1.29%  "
$ws.Range("E37").Value = "  +5.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.90"
$ws.Range("E38").Value = "  +2.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.040.79"
$ws.Range("E39").Value = "  +2.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.98"
$ws.Range("E40").Value = "  +3.31%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.650"
$ws.Range("E42").Value = "  +3.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.282.72"
$ws.Range("E43").Value = "  +7.95%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.69"
$ws.Range("E44").Value = "  +5.43%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.01"
$ws.Range("E45").Value = "  +1.72%  "
$ws.Range("E46").Value = "  +3.86%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0242"
$ws.Range("E47").Value = "  +6.13%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.96"
$ws.Range("E48").Value = "  +14.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.86"
$ws.Range("E49").Value = "  +5.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.38"
$ws.Range("E50").Value = "  +3.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0874"
$ws.Range("E51").Value = "  +6.55%  "
